# feat: add 2022-Q1 data
#
# The workbook currently has sheets: 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计.
# This script:
#   1. Inserts a new "2022-Q1" sheet (same per-fund layout as the other
#      quarterly sheets) right before the "总计" summary sheet.
#   2. Populates it with the Q1-2022 fund holdings.
#   3. Prepends a new "2022-Q1" row to the "总计" summary sheet, shifting the
#      existing rows down and renumbering the index column.
#
# NOTE: worksheet handles returned before a sheet is inserted/moved can go
# stale (they track sheet *position*, not identity), so sheets are always
# re-fetched by name with Worksheets.Item(...) right before they're used.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new quarter sheet right before "总计" ---------------------
$beforeSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q1"

# Re-fetch fresh handles by name now that the sheet collection has changed.
$q1 = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")

# --- 2. Fill in the 2022-Q1 per-fund data -----------------------------------
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1Data = @(
    @(0, "539003", "建信富时100指数（QDII）人民币A",   "0.71", "92.86", "5.73", "0.0407", 4),
    @(1, "008707", "建信富时100指数（QDII）美元现汇A", "0.71", "92.86", "5.73", "0.0407", 4),
    @(2, "008706", "建信富时100指数（QDII）人民币C",   "0.20", "92.86", "5.73", "0.0115", 4),
    @(3, "008708", "建信富时100指数（QDII）美元现汇C", "0.20", "92.86", "5.73", "0.0115", 4)
)

foreach ($row in $q1Data) {
    $r = [int]$row[0] + 2
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
}

# --- 3. Prepend the 2022-Q1 row to the "总计" sheet -------------------------
# Existing rows (date, count, market value) all shift down by one, and the
# leading index column (A) is renumbered 0..4.
$totalData = @(
    @("2022-Q1", 4, 0.1),
    @("2021-Q4", 4, 0.43),
    @("2021-Q3", 4, 0.14),
    @("2021-Q2", 4, 0.18),
    @("2021-Q1", 4, 0.17)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $totalData[$i][0]
    $totalSheet.Cells.Item($r, 3).Value = $totalData[$i][1]
    $totalSheet.Cells.Item($r, 4).Value = $totalData[$i][2]
}
